# Apply the "effort" log update:
# - Add three new shared-string rows (new log entries)
# - Append three new data rows (42, 43, 44) with date / effort / comment
# - Update the view (scrolled down, new active cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell that already carries the date-column style (custom format
# "ddd dd/mm/yyyy"). Copying its format (instead of re-typing the format
# string) keeps the original style index / numFmt instead of minting a
# duplicate one.
$refCell = $ws.Range("A41")

# New log entries to append
$newRows = @(
    @{ Row = 42; Date = 41228; Effort = 2.5;  Comment = "Test case tc08 is working well, a fix of rtos.c was required. All test cases rerun. Manual updated" },
    @{ Row = 43; Date = 41229; Effort = 3.25; Comment = "Preparation of release, new test case tc09" },
    @{ Row = 44; Date = 41232; Effort = 1.25; Comment = "Makefile: Workaround for 12 Bit Branch distance problem with core.a" }
)

foreach ($r in $newRows) {
    $rowIdx = $r.Row

    $cellA = $ws.Cells.Item($rowIdx, 1)
    $refCell.Copy()
    $cellA.PasteSpecial(-4122)  # xlPasteFormats
    $cellA.Value = $r.Date

    $cellB = $ws.Cells.Item($rowIdx, 2)
    $cellB.Value = $r.Effort

    $cellD = $ws.Cells.Item($rowIdx, 4)
    $cellD.Value = $r.Comment
}

# Update the visible view: scroll so row 38 is at top and select A44
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 38
$ws.Range("A44").Select() | Out-Null
